# resumo_campanha_63.xlsx - "Sistmea de chamados upgrade"
# Updates the VENDEDOR sales figures (and downstream totals / summary
# tables) on the "RESUMO CAMPANHA 63" sheet to the new period numbers,
# and re-sorts the sales-rep names alphabetically (Alan, Ana, Sibely,
# Sttefani, Suzana) in every table that lists them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- VENDA EM BOLETOS - MATRIZ (rows 8-12) + TOTAL (row 13) ---------------
# Row 8: Alan
$ws.Range("B8").Value = "Alan"
$ws.Range("C8").Value = 280
$ws.Range("D8").Value = 280
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

# Row 9: Ana
$ws.Range("B9").Value = "Ana"
$ws.Range("C9").Value = 6879
$ws.Range("D9").Value = 4979
$ws.Range("E9").Value = 1900
$ws.Range("F9").Value = 0.27620293647332
$ws.Range("G9").Value = 19

# Row 10: Sibely
$ws.Range("B10").Value = "Sibely"
$ws.Range("C10").Value = 4080
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 4080
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 40.8

# Row 11: Sttefani - unchanged

# Row 12: Suzana
$ws.Range("B12").Value = "Suzana"
$ws.Range("C12").Value = 8750
$ws.Range("D12").Value = 7700
$ws.Range("E12").Value = 1050
$ws.Range("F12").Value = 0.12
$ws.Range("G12").Value = 10.5

# Row 13: TOTAL
$ws.Range("C13").Value = 31443
$ws.Range("D13").Value = 19669
$ws.Range("E13").Value = 11774
$ws.Range("F13").Value = 0.37445536367395
$ws.Range("G13").Value = 117.74

# --- VENDA EM BOLETOS - DAILY (row 21) + TOTAL (row 22) -------------------
$ws.Range("D21").Value = 2600
$ws.Range("E21").Value = 2315
$ws.Range("F21").Value = 0.47100712105799
$ws.Range("G21").Value = 23.15

$ws.Range("D22").Value = 2600
$ws.Range("E22").Value = 2315
$ws.Range("F22").Value = 0.47100712105799
$ws.Range("G22").Value = 23.15

# --- RESUMO GERAL VENDAS EM BOLETOS (rows 25, 27, 28) ----------------------
$ws.Range("C25").Value = 31443
$ws.Range("D25").Value = 19669
$ws.Range("E25").Value = 11774
$ws.Range("F25").Value = 0.37445536367395
$ws.Range("G25").Value = 117.74

$ws.Range("D27").Value = 2600
$ws.Range("E27").Value = 2315
$ws.Range("F27").Value = 0.47100712105799
$ws.Range("G27").Value = 23.15

$ws.Range("C28").Value = 36678
$ws.Range("D28").Value = 22269
$ws.Range("E28").Value = 14409
$ws.Range("F28").Value = 0.39285130050712
$ws.Range("G28").Value = 144.09

# --- OTICAS breakdown block (rows 32-36, 38, 40, 41) -----------------------
$ws.Range("C32").Value = "Alan"
$ws.Range("D32").Value = 280
$ws.Range("E32").Value = 280
$ws.Range("F32").Value = 0

$ws.Range("C33").Value = "Ana"
$ws.Range("D33").Value = 6879
$ws.Range("E33").Value = 4979
$ws.Range("F33").Value = 1900

$ws.Range("C34").Value = "Sibely"
$ws.Range("D34").Value = 4080
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 4080

# Row 35: Sttefani - unchanged

$ws.Range("C36").Value = "Suzana"
$ws.Range("D36").Value = 8750
$ws.Range("E36").Value = 7700
$ws.Range("F36").Value = 1050

$ws.Range("E40").Value = 2600
$ws.Range("F40").Value = 2315

$ws.Range("D41").Value = 36678
$ws.Range("E41").Value = 22269
$ws.Range("F41").Value = 14409

# CONVERTIDO / NÃO CONVERTIDO percentages
$ws.Range("E42").Value = 0.60714869949288
$ws.Range("E43").Value = 0.39285130050712
